# Update line power-flow results for the 380 kV case
# (Code/Results/Cases/Case_3_252/res_line/pl_mw.xlsx)
# Rows 2-25 (index 0-23); columns B,C,D,F,G,I,J,K updated, others unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @(2, "B", 0.5403633160001675),
    @(2, "C", 0.2294126817136259),
    @(2, "D", 0.08405762848991571),
    @(2, "F", 3.418410032993052),
    @(2, "G", 0.002562522143272855),
    @(2, "I", 1.793102721377323),
    @(2, "J", 0.3562383099905162),
    @(2, "K", 0.7682098671597259),
    @(3, "B", 0.5142995554442678),
    @(3, "C", 0.2196148699795515),
    @(3, "D", 0.08275542837323968),
    @(3, "F", 3.365981028969955),
    @(3, "G", 0.002566957033344344),
    @(3, "I", 1.768429233866755),
    @(3, "J", 0.3456689397605714),
    @(3, "K", 0.7325232995724775),
    @(4, "B", 0.4987180274170555),
    @(4, "C", 0.2137535101394121),
    @(4, "D", 0.08203794720451185),
    @(4, "F", 3.335049297300969),
    @(4, "G", 0.002569822436095227),
    @(4, "I", 1.753800959148847),
    @(4, "J", 0.3393993371741573),
    @(4, "K", 0.7111861989267538),
    @(5, "B", 0.4924743911035137),
    @(5, "C", 0.2114036467890799),
    @(5, "D", 0.08176622478332263),
    @(5, "F", 3.322760078369711),
    @(5, "G", 0.002571026029129726),
    @(5, "I", 1.747970301765008),
    @(5, "J", 0.3368995280341096),
    @(5, "K", 0.7026353740993443),
    @(6, "B", 0.4914440377148424),
    @(6, "C", 0.211015785521198),
    @(6, "D", 0.08172235368840575),
    @(6, "F", 3.320738499371359),
    @(6, "G", 0.002571228057680952),
    @(6, "I", 1.747009988539972),
    @(6, "J", 0.3364877583382651),
    @(6, "K", 0.701224218883624),
    @(7, "B", 0.4986333946210664),
    @(7, "C", 0.213721662638136),
    @(7, "D", 0.08203419899599851),
    @(7, "F", 3.334882283756627),
    @(7, "G", 0.002569838522635991),
    @(7, "I", 1.753721797366723),
    @(7, "J", 0.3393654010487808),
    @(7, "K", 0.7110702959350022),
    @(8, "B", 0.5312889391630051),
    @(8, "C", 0.2260022149123415),
    @(8, "D", 0.08359160512893737),
    @(8, "F", 3.40007033929173),
    @(8, "G", 0.00256402181874682),
    @(8, "I", 1.784486718412424),
    @(8, "J", 0.3525481832290041),
    @(8, "K", 0.7557857019419885),
    @(9, "B", 0.5986810493199073),
    @(9, "C", 0.2513198730108002),
    @(9, "D", 0.08729651905011337),
    @(9, "F", 3.537959364043047),
    @(9, "G", 0.002553739372163563),
    @(9, "I", 1.848984991112914),
    @(9, "J", 0.3801574325953254),
    @(9, "K", 0.8480501159391736),
    @(10, "B", 0.6502576733289231),
    @(10, "C", 0.2706892992007965),
    @(10, "D", 0.09041542087480536),
    @(10, "F", 3.645491386753747),
    @(10, "G", 0.002546862435469042),
    @(10, "I", 1.898962735687391),
    @(10, "J", 0.4015328741461133),
    @(10, "K", 0.9186621139127737),
    @(11, "B", 0.6741738695352808),
    @(11, "C", 0.2796713322285029),
    @(11, "D", 0.0919205681261559),
    @(11, "F", 3.695783593004535),
    @(11, "G", 0.002543879425119745),
    @(11, "I", 1.922272585356311),
    @(11, "J", 0.4114983735017574),
    @(11, "K", 0.9514070552000078),
    @(12, "B", 0.6832958054264395),
    @(12, "C", 0.2830973940926071),
    @(12, "D", 0.09250294307778972),
    @(12, "F", 3.715027189942873),
    @(12, "G", 0.002542770611540796),
    @(12, "I", 1.931182842852209),
    @(12, "J", 0.4153071114753288),
    @(12, "K", 0.963896830646064),
    @(13, "B", 0.6813283225964426),
    @(13, "C", 0.2823584258589733),
    @(13, "D", 0.09237696645288906),
    @(13, "F", 3.710873864969329),
    @(13, "G", 0.002543008491482815),
    @(13, "I", 1.929260142790952),
    @(13, "J", 0.4144852698669581),
    @(13, "K", 0.9612029263642512),
    @(14, "B", 0.6749230252810605),
    @(14, "C", 0.2799526989822141),
    @(14, "D", 0.09196823182051617),
    @(14, "F", 3.697362779208845),
    @(14, "G", 0.002543787786143428),
    @(14, "I", 1.923003963777347),
    @(14, "J", 0.4118110172997262),
    @(14, "K", 0.9524327911456112),
    @(15, "B", 0.6710081161056678),
    @(15, "C", 0.2784823524781359),
    @(15, "D", 0.09171948578780587),
    @(15, "F", 3.689112806132329),
    @(15, "G", 0.002544267830702575),
    @(15, "I", 1.919182744859157),
    @(15, "J", 0.4101775288523868),
    @(15, "K", 0.947072558736437),
    @(16, "B", 0.6487038361024418),
    @(16, "C", 0.2701057560375943),
    @(16, "D", 0.0903187927145126),
    @(16, "F", 3.642232455445338),
    @(16, "G", 0.002547060297254788),
    @(16, "I", 1.897451005024394),
    @(16, "J", 0.4008864914073911),
    @(16, "K", 0.9165347258729071),
    @(17, "B", 0.6351371945945914),
    @(17, "C", 0.2650108700676697),
    @(17, "D", 0.0894816214777876),
    @(17, "F", 3.61382602857401),
    @(17, "G", 0.00254881053079501),
    @(17, "I", 1.884266968601295),
    @(17, "J", 0.3952488327885391),
    @(17, "K", 0.897960622251361),
    @(18, "B", 0.6273767192594164),
    @(18, "C", 0.2620964890664652),
    @(18, "D", 0.08900823011617831),
    @(18, "F", 3.597616781726515),
    @(18, "G", 0.002549830905870624),
    @(18, "I", 1.87673791439272),
    @(18, "J", 0.3920289349961052),
    @(18, "K", 0.8873359241087257),
    @(19, "B", 0.6247564848093532),
    @(19, "C", 0.2611124824087199),
    @(19, "D", 0.0888493440832292),
    @(19, "F", 3.592150792067514),
    @(19, "G", 0.002550178741489517),
    @(19, "I", 1.874197969348486),
    @(19, "J", 0.3909426296285545),
    @(19, "K", 0.8837486436488007),
    @(20, "B", 0.6365769666418544),
    @(20, "C", 0.2655515654996634),
    @(20, "D", 0.08956989869209053),
    @(20, "F", 3.61683654200408),
    @(20, "G", 0.002548622799659998),
    @(20, "I", 1.885664832434841),
    @(20, "J", 0.3958466163606005),
    @(20, "K", 0.8999317962104669),
    @(21, "B", 0.6768026404566569),
    @(21, "C", 0.280658645645758),
    @(21, "D", 0.09208795038445317),
    @(21, "F", 3.701325900916743),
    @(21, "G", 0.002543558325214231),
    @(21, "I", 1.924839288403419),
    @(21, "J", 0.4125955575645008),
    @(21, "K", 0.9550063475446962),
    @(22, "B", 0.7034736695147217),
    @(22, "C", 0.2906764179596735),
    @(22, "D", 0.09380596564911059),
    @(22, "F", 3.757705272303582),
    @(22, "G", 0.002540369513800956),
    @(22, "I", 1.950928050066707),
    @(22, "J", 0.4237462051919039),
    @(22, "K", 0.9915253826034416),
    @(23, "B", 0.6892039065488405),
    @(23, "C", 0.2853164618238395),
    @(23, "D", 0.09288241315186951),
    @(23, "F", 3.727507918851927),
    @(23, "G", 0.00254206039755825),
    @(23, "I", 1.936959293456539),
    @(23, "J", 0.4177761181758655),
    @(23, "K", 0.9719863589715487),
    @(24, "B", 0.6359259237089532),
    @(24, "C", 0.2653070712083831),
    @(24, "D", 0.08952996392605428),
    @(24, "F", 3.615475108993706),
    @(24, "G", 0.002548707628772973),
    @(24, "I", 1.885032700552117),
    @(24, "J", 0.395576292187144),
    @(24, "K", 0.8990404609082816),
    @(25, "B", 0.5800886751059124),
    @(25, "C", 0.2443369323956688),
    @(25, "D", 0.08622455646903404),
    @(25, "F", 3.499570460167689),
    @(25, "G", 0.002556401502727978),
    @(25, "I", 1.831085121459367),
    @(25, "J", 0.3724982288705547),
    @(25, "K", 0.8225965427106416)
)

foreach ($entry in $updates) {
    $r = $entry[0]
    $col = $entry[1]
    $val = $entry[2]
    $ws.Range("$col$r").Value = $val
}
